# Update the workbook to match the target revision:
#  - G21/H21 values updated (1.612 -> 1.614, 77 -> 78)
#  - Row 37 (Squilla mantis / SQUIMAN duplicate entry) removed entirely,
#    shifting all subsequent rows up by one and shrinking the used range
#    from A1:K41 to A1:K40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the weight/number for row 21 (Squilla mantis, first 2-RAP occurrence)
$ws.Range("G21").Value = 1.614
$ws.Range("H21").Value = 78

# Remove row 37 entirely (duplicate Squilla mantis entry), shifting rows
# 38-41 up to become rows 37-40.
$ws.Rows(37).Delete()
